$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row height tweaks (platform re-save rounding: 38.1->38, 17.1->17, 15.95->16)
$ws.Rows.Item(2).RowHeight = 38
$ws.Range("A3:A28").EntireRow.RowHeight = 17
$ws.Range("A29:A33").EntireRow.RowHeight = 16

# Row 4: replace "Example 1" / date / hours / description
$ws.Range("A4").Value = "Downloaded the live version of the app"
$ws.Range("B4").Value = (Get-Date -Year 2019 -Month 1 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = "Explored the already made app to inform my app"

# Row 5: clear "Example 2" subject text (keep date/hours the same), update description text (same text, just reindexed in sharedStrings so no visible change)
$ws.Range("A5").Value = ""
$ws.Range("D5").Value = "Had some issues with…"

# Sheet view changes: zoom and selection
$ws.Application.ActiveWindow.Zoom = 175
$ws.Range("A5").Select()
